$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and report date range) ---
$ws.Range("A8").Value = "Volume 32   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# --- Cells changing FROM numeric TO shared-string text "0" (keep style 13) ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C28").PasteSpecial(-4122)

# --- Cell changing FROM numeric TO shared-string text "***.*" (keep style 13) ---
$ws.Range("E22").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# --- Cells changing FROM shared-string text TO numeric (style must switch to 14/15) ---
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 2
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E29").Value = -50
$ws.Range("E30").Value = 0
$ws.Range("H14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E30").PasteSpecial(-4122)

# --- Remaining plain numeric value updates (style unchanged) ---
$ws.Range("M14").Value = -16.666666666666
$ws.Range("E15").Value = -100
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 105.882352941176
$ws.Range("L15").Value = 29.629629629629
$ws.Range("M15").Value = 94.444444444444
$ws.Range("C16").Value = 13
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 85.714285714285
$ws.Range("G16").Value = 39
$ws.Range("H16").Value = 23.076923076923
$ws.Range("I16").Value = 263
$ws.Range("J16").Value = 274
$ws.Range("K16").Value = -4.014598540145
$ws.Range("L16").Value = -10.847457627118
$ws.Range("M16").Value = -0.378787878787
$ws.Range("N16").Value = -75.937785910338
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = -41.176470588235
$ws.Range("F17").Value = 74
$ws.Range("G17").Value = 58
$ws.Range("H17").Value = 27.586206896551
$ws.Range("I17").Value = 428
$ws.Range("J17").Value = 402
$ws.Range("K17").Value = 6.467661691542
$ws.Range("L17").Value = -9.129511677282
$ws.Range("M17").Value = 56.204379562043
$ws.Range("N17").Value = -14.910536779324
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 20
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 153
$ws.Range("J18").Value = 156
$ws.Range("K18").Value = -1.923076923076
$ws.Range("L18").Value = -22.33502538071
$ws.Range("M18").Value = -25.365853658536
$ws.Range("N18").Value = -84.451219512195
$ws.Range("C19").Value = 20
$ws.Range("E19").Value = -20
$ws.Range("F19").Value = 93
$ws.Range("G19").Value = 96
$ws.Range("H19").Value = -3.125
$ws.Range("I19").Value = 625
$ws.Range("J19").Value = 548
$ws.Range("K19").Value = 14.05109489051
$ws.Range("L19").Value = 19.502868068833
$ws.Range("M19").Value = 102.922077922078
$ws.Range("N19").Value = 55.860349127182
$ws.Range("C20").Value = 12
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 20
$ws.Range("G20").Value = 51
$ws.Range("H20").Value = -15.686274509803
$ws.Range("I20").Value = 280
$ws.Range("J20").Value = 264
$ws.Range("K20").Value = 6.060606060606
$ws.Range("L20").Value = -25.333333333333
$ws.Range("M20").Value = 113.740458015267
$ws.Range("N20").Value = -72.468043264503
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 65
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 282
$ws.Range("G21").Value = 271
$ws.Range("H21").Value = 4.059040590405
$ws.Range("I21").Value = 1789
$ws.Range("J21").Value = 1667
$ws.Range("K21").Value = 7.318536292741
$ws.Range("L21").Value = -5.593667546174
$ws.Range("M21").Value = 48.341625207296
$ws.Range("N21").Value = -56.130456105934
$ws.Range("M22").Value = -45.454545454545
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -77.777777777777
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = -36.666666666666
$ws.Range("I23").Value = 130
$ws.Range("J23").Value = 151
$ws.Range("K23").Value = -13.907284768211
$ws.Range("L23").Value = -23.076923076923
$ws.Range("M23").Value = 6.55737704918
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 33
$ws.Range("E24").Value = -36.363636363636
$ws.Range("F24").Value = 144
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = 18.032786885245
$ws.Range("I24").Value = 1160
$ws.Range("J24").Value = 974
$ws.Range("K24").Value = 19.096509240246
$ws.Range("L24").Value = 2.654867256637
$ws.Range("M24").Value = 48.148148148148
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -71.428571428571
$ws.Range("F25").Value = 47
$ws.Range("H25").Value = -22.950819672131
$ws.Range("I25").Value = 356
$ws.Range("J25").Value = 394
$ws.Range("K25").Value = -9.644670050761
$ws.Range("L25").Value = -24.735729386892
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 23
$ws.Range("E26").Value = 8.695652173913
$ws.Range("F26").Value = 97
$ws.Range("G26").Value = 83
$ws.Range("H26").Value = 16.867469879518
$ws.Range("I26").Value = 607
$ws.Range("J26").Value = 658
$ws.Range("K26").Value = -7.750759878419
$ws.Range("L26").Value = -5.598755832037
$ws.Range("M26").Value = -31.950672645739
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = 35.483870967741
$ws.Range("L27").Value = 16.666666666666
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 64
$ws.Range("K28").Value = -39.0625
$ws.Range("L28").Value = -35
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 17
$ws.Range("J29").Value = 16
$ws.Range("K29").Value = 6.25
$ws.Range("L29").Value = -5.555555555555
$ws.Range("M29").Value = -34.615384615384
$ws.Range("N29").Value = -78.75
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 150
$ws.Range("I30").Value = 16
$ws.Range("J30").Value = 14
$ws.Range("K30").Value = 14.285714285714
$ws.Range("L30").Value = 6.666666666666
$ws.Range("M30").Value = -30.434782608695
$ws.Range("N30").Value = -78.666666666666
